$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.809.91'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.639.36'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.06'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("E6").Value = '  -0.65%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.258'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.69'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0791'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.26'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.865.76'
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = '1.639.21'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.07'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '25.844.89'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("E19").Value = '  -0.18%  '
$ws.Range("E20").Value = '  +2.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.71'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.87%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  +4.22%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.27%  '
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.94'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.33'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("E33").Value = '  -0.28%  '
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("E35").Value = '  -0.29%  '
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").Value = '1.132.51'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.54'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.547'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.18%  '
$ws.Range("E40").Value = '  +0.54%  '
$ws.Range("E41").Value = '  -0.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.59'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.71'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.805'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("D45").Value = '1.775.52'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("E46").Value = '  +3.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.35'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.13%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.417'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.44'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.33%  '
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.48'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.06%  '
